# Weekly refresh of the Fruta/Hortaliza (Damasco) dataset:
# a new daily record is inserted at the top of the data block (row 10),
# pushing every existing record down by one row (old row 24 -> new row 25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 10, shifting rows 10-24
# down to 11-25 (carrying their formatting, e.g. the date style on column D).
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new record.
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C10").Value = "Metropolitana"
$ws.Range("D10").Value = 44168
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100103
$ws.Range("H10").Value = "Frutos de hueso (carozo)"
$ws.Range("I10").Value = 100103003
$ws.Range("J10").Value = "Damasco"
$ws.Range("K10").Value = "Dina"
$ws.Range("L10").Value = "Especial"
$ws.Range("M10").Value = 40
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 14000
$ws.Range("P10").Value = 14000
$ws.Range("Q10").Value = "$/bandeja 10 kilos"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 1400
$ws.Range("T10").Value = 10
